$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 836.76
$wsSummary.Range("E2").Value = 9163.24
$wsSummary.Range("F2").Value = 849.4

$wsSummary.Range("A3").Value = 561.21
$wsSummary.Range("E3").Value = 510.25
$wsSummary.Range("F3").Value = 38.32

$wsSummary.Range("A5").Value = 200
$wsSummary.Range("B5").Value = 100
$wsSummary.Range("E5").Value = 100
$wsSummary.Range("F5").Value = 100

$wsSummary.Range("F5").Select()

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

$wsSchedule.Range("F5").Value = 849.4
$wsSchedule.Range("G5").Value = 8313.84
$wsSchedule.Range("H5").Value = 38.32
$wsSchedule.Range("J5").Value = 100
$wsSchedule.Range("K5").Value = 987.72
$wsSchedule.Range("P5").Value = 987.72

$wsSchedule.Range("F6").Value = 794.33
$wsSchedule.Range("G6").Value = 7519.51
$wsSchedule.Range("H6").Value = 93.39

$wsSchedule.Range("F7").Value = 813.55
$wsSchedule.Range("G7").Value = 6705.96
$wsSchedule.Range("H7").Value = 74.17

$wsSchedule.Range("F8").Value = 819.37
$wsSchedule.Range("G8").Value = 5886.59
$wsSchedule.Range("H8").Value = 68.35

$wsSchedule.Range("F9").Value = 829.66
$wsSchedule.Range("G9").Value = 5056.93
$wsSchedule.Range("H9").Value = 58.06

$wsSchedule.Range("F10").Value = 836.18
$wsSchedule.Range("G10").Value = 4220.75
$wsSchedule.Range("H10").Value = 51.54

$wsSchedule.Range("F11").Value = 844.7
$wsSchedule.Range("G11").Value = 3376.05
$wsSchedule.Range("H11").Value = 43.02

$wsSchedule.Range("F12").Value = 854.42
$wsSchedule.Range("G12").Value = 2521.63
$wsSchedule.Range("H12").Value = 33.3

$wsSchedule.Range("F13").Value = 862.02
$wsSchedule.Range("G13").Value = 1659.61
$wsSchedule.Range("H13").Value = 25.7

$wsSchedule.Range("F14").Value = 871.35
$wsSchedule.Range("G14").Value = 788.26
$wsSchedule.Range("H14").Value = 16.37

$wsSchedule.Range("F15").Value = 788.26
$wsSchedule.Range("H15").Value = 8.03
$wsSchedule.Range("K15").Value = 796.29
$wsSchedule.Range("P15").Value = 796.29

# Drop the now-unused "Over Due" column (O) entries and the stray P2 cell
# so the underlying cells disappear rather than linger as empty values.
$wsSchedule.Range("P2").Clear()
$wsSchedule.Range("O3").Clear()
$wsSchedule.Range("O4").Clear()
$wsSchedule.Range("O5:O15").Clear()

$wsSchedule.Range("G15").Select()

# ---------------------------------------------------------------------------
# Transactions sheet - insert a new disbursement transaction row at the top
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Rows("2:2").Insert()

# Match formatting of the surrounding rows before filling in values:
# columns A-I follow the "disbursement" row template, J-L follow the
# "repayment" row template (so J/K/L pick up the right number styles).
$wsTransactions.Range("A4:I4").Copy()
$wsTransactions.Range("A2:I2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsTransactions.Range("J3:L3").Copy()
$wsTransactions.Range("J2:L2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsTransactions.Range("A2").Value = 1671
$wsTransactions.Range("B2").Value = "Head Office"
$wsTransactions.Range("C2").Value = 42064
$wsTransactions.Range("D2").Value = "Disbursement"
$wsTransactions.Range("E2").Value = 5000
$wsTransactions.Range("F2").Value = 0
$wsTransactions.Range("G2").Value = 0
$wsTransactions.Range("H2").Value = 0
$wsTransactions.Range("I2").Value = 0
$wsTransactions.Range("J2").Value = 9163.24

$wsTransactions.Range("A3").Value = 640
$wsTransactions.Range("A4").Value = 632

$wsTransactions.Activate()
$wsTransactions.Range("J3").Select()
